# Adapt column header formatting to respective input file names:
#   "<...>_old" -> "<...>_FV2404"
#   "<...>_new" -> "<...>_FV2410"
# and turn the header+data range into an Excel Table (with AutoFilter),
# freezing the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base names shared by the "_old"/"_new" (now "_FV2404"/"_FV2410") column pairs,
# in left-to-right sheet order.
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J: "<name>_old" -> "<name>_FV2404"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2404"
}

# Column K ("diff") is unchanged.
$ws.Cells.Item(1, 11).Value = "diff"

# Columns L-U: "<name>_new" -> "<name>_FV2410"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2410"
}

# Convert the used range (header row 1 + data rows 2-55, columns A-U) into a
# native Excel Table so the workbook gets an autofilter + table definition.
$dataRange = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (pane split below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
